# dlgCamperInsurance.xlsx — add the per-section XPath locator columns (C1:G1)
# to the header row and move the sheet selection onto them.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$ws.Range("C1").Value = "//*[@id='insurance-form']/div/section[1]"
$ws.Range("D1").Value = "//*[@id='insurance-form']/div/section[2]"
$ws.Range("E1").Value = "//*[@id='insurance-form']/div/section[3]"
$ws.Range("F1").Value = "//*[@id='insurance-form']/div/section[4]"
$ws.Range("G1").Value = "//*[@id='insurance-form']/div/section[5]"

# Match the author's new selection (was A13) onto the newly-added cells.
$ws.Range("C1:G1").Select()
